$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string text updates (report volume/date header text) ---
$ws.Range("A8").Value = "Volume 30   Number  6"
$ws.Range("C9").Value = "Report Covering the Week  2/6/2023  Through  2/12/2023"

# --- Crime statistics table updates (rows 16-30) ---
# Row 16
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = 300
$ws.Range("I16").Value = 13
$ws.Range("J16").Value = 16
$ws.Range("K16").Value = -18.75
$ws.Range("L16").Value = 225
$ws.Range("M16").Value = -38.095238095238
$ws.Range("N16").Value = -81.159420289855
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "0"
$ws.Range("A16").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("D16").Value = 1
$ws.Range("F16").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").Value = -100
$ws.Range("L16").Copy()
$ws.Range("E16").PasteSpecial(-4122)

# Row 17
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 200
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = -16.666666666666
$ws.Range("I17").Value = 13
$ws.Range("J17").Value = 16
$ws.Range("K17").Value = -18.75
$ws.Range("L17").Value = 62.5
$ws.Range("M17").Value = 85.714285714285
$ws.Range("N17").Value = -48

# Row 18
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -80
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = -6.666666666666
$ws.Range("I18").Value = 28
$ws.Range("J18").Value = 33
$ws.Range("K18").Value = -15.151515151515
$ws.Range("L18").Value = -3.448275862068
$ws.Range("M18").Value = 21.739130434782
$ws.Range("N18").Value = -81.699346405228

# Row 19
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 22.222222222222
$ws.Range("F19").Value = 51
$ws.Range("G19").Value = 32
$ws.Range("H19").Value = 59.375
$ws.Range("I19").Value = 69
$ws.Range("J19").Value = 48
$ws.Range("K19").Value = 43.75
$ws.Range("L19").Value = 60.465116279069
$ws.Range("M19").Value = 146.428571428571
$ws.Range("N19").Value = 53.333333333333

# Row 20
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 133.333333333333
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = 18.181818181818
$ws.Range("I20").Value = 16
$ws.Range("J20").Value = 12
$ws.Range("K20").Value = 33.333333333333
$ws.Range("L20").Value = 128.571428571429
$ws.Range("M20").Value = 77.777777777777
$ws.Range("N20").Value = -83.157894736842

# Row 21
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = 15.789473684210
$ws.Range("F21").Value = 96
$ws.Range("G21").Value = 72
$ws.Range("H21").Value = 33.333333333333
$ws.Range("I21").Value = 139
$ws.Range("J21").Value = 125
$ws.Range("K21").Value = 11.2
$ws.Range("L21").Value = 51.086956521739
$ws.Range("M21").Value = 57.954545454545
$ws.Range("N21").Value = -64.26735218509

# Row 23
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 4
$ws.Range("J23").Value = 2
$ws.Range("K23").Value = 100
$ws.Range("L23").Value = 100

# Row 24
$ws.Range("C24").Value = 25
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = 56.25
$ws.Range("F24").Value = 72
$ws.Range("G24").Value = 64
$ws.Range("H24").Value = 12.5
$ws.Range("I24").Value = 101
$ws.Range("J24").Value = 104
$ws.Range("K24").Value = -2.884615384615
$ws.Range("L24").Value = 42.253521126760
$ws.Range("M24").Value = 102

# Row 25
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 13
$ws.Range("G25").Value = 18
$ws.Range("H25").Value = -27.777777777777
$ws.Range("I25").Value = 27
$ws.Range("J25").Value = 27
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 12.5
$ws.Range("M25").Value = 35

# Row 27
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 133.333333333333
$ws.Range("I27").Value = 9
$ws.Range("J27").Value = 3
$ws.Range("K27").Value = 200
$ws.Range("L27").Value = 125
$ws.Range("C27").Value = 2
$ws.Range("F16").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D27").Value = 1
$ws.Range("F16").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = 100
$ws.Range("L16").Copy()
$ws.Range("E27").PasteSpecial(-4122)

# Row 30
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "0"
$ws.Range("A16").Copy()
$ws.Range("C30").PasteSpecial(-4122)

$excel.CutCopyMode = $false
Write-Output "Edits applied successfully"